# Fruta / hortaliza, semanal
# Insert 2 new rows at row 979 (pushing existing data for rows 979-1068 down to 981-1070)
# and populate the two new rows with the latest weekly price observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 979; existing rows shift down by one each time.
$ws.Rows.Item(979).Insert()
$ws.Rows.Item(979).Insert()

# --- New row 979 ---
$ws.Cells.Item(979, 1).Value = 9
$ws.Cells.Item(979, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(979, 3).Value = "Metropolitana"
$ws.Cells.Item(979, 4).Value = 44946
$ws.Cells.Item(979, 5).Value = 13
$ws.Cells.Item(979, 6).Value = "Fruta"
$ws.Cells.Item(979, 7).Value = 100102
$ws.Cells.Item(979, 8).Value = "Cítricos"
$ws.Cells.Item(979, 9).Value = 100102005
$ws.Cells.Item(979, 10).Value = "Naranja"
$ws.Cells.Item(979, 11).Value = "Valencia"
$ws.Cells.Item(979, 12).Value = "Primera"
$ws.Cells.Item(979, 13).Value = 470
$ws.Cells.Item(979, 14).Value = 10500
$ws.Cells.Item(979, 15).Value = 11000
$ws.Cells.Item(979, 16).Value = 10734
$ws.Cells.Item(979, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(979, 18).Value = "Cabildo"
$ws.Cells.Item(979, 19).Value = 716
$ws.Cells.Item(979, 20).Value = 15

# --- New row 980 ---
$ws.Cells.Item(980, 1).Value = 9
$ws.Cells.Item(980, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(980, 3).Value = "Metropolitana"
$ws.Cells.Item(980, 4).Value = 44946
$ws.Cells.Item(980, 5).Value = 13
$ws.Cells.Item(980, 6).Value = "Fruta"
$ws.Cells.Item(980, 7).Value = 100102
$ws.Cells.Item(980, 8).Value = "Cítricos"
$ws.Cells.Item(980, 9).Value = 100102005
$ws.Cells.Item(980, 10).Value = "Naranja"
$ws.Cells.Item(980, 11).Value = "Valencia"
$ws.Cells.Item(980, 12).Value = "Primera"
$ws.Cells.Item(980, 13).Value = 450
$ws.Cells.Item(980, 14).Value = 11500
$ws.Cells.Item(980, 15).Value = 12000
$ws.Cells.Item(980, 16).Value = 11722
$ws.Cells.Item(980, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(980, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(980, 19).Value = 651
$ws.Cells.Item(980, 20).Value = 18
